$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = "A-Failed"
$ws.Range("E3").Select()
